# "removed false start data": the first two logged rows (2 and 3) were an
# erroneous early run and should be dropped, shifting the remaining two
# data rows (old rows 4 and 5) up into rows 2 and 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:3").Select()
$ws.Rows("2:3").Delete()
